$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last existing header cell (H1) onto the
# two new header cells so they pick up the same bold/border/alignment style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add the new data values (rows 2 and 3)
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
